$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear cells that no longer exist in the updated layout (row 7, columns A-H)
$ws.Range("A7:H7").Clear()

# Row 1
$ws.Range("A1").Value = "negative"
$ws.Range("J1").Value = "positive"

# Row 2
$ws.Range("A2").Value = "name"
$ws.Range("B2").Value = "anchor score"
$ws.Range("C2").Value = "type occurences"
$ws.Range("D2").Value = "total occurences"
$ws.Range("E2").Value = "+%"
$ws.Range("F2").Value = "-%"
$ws.Range("G2").Value = "both"
$ws.Range("H2").Value = "normal"
$ws.Range("J2").Value = "name"
$ws.Range("K2").Value = "anchor score"
$ws.Range("L2").Value = "type occurences"
$ws.Range("M2").Value = "total occurences"
$ws.Range("N2").Value = "+%"
$ws.Range("O2").Value = "-%"
$ws.Range("P2").Value = "both"
$ws.Range("Q2").Value = "normal"

# Row 3
$ws.Range("A3").Value = "crude"
$ws.Range("B3").Value = 0.8235294117647058
$ws.Range("C3").Value = 28
$ws.Range("D3").Value = 28
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = $false
$ws.Range("H3").Value = 6
$ws.Range("J3").Value = "interesting"
$ws.Range("K3").Value = 0.9696969696969697
$ws.Range("L3").Value = 32
$ws.Range("M3").Value = 32
$ws.Range("N3").Value = 1
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = $false
$ws.Range("Q3").Value = 1

# Row 4
$ws.Range("A4").Value = "crisis"
$ws.Range("B4").Value = 0.5958904109589042
$ws.Range("C4").Value = 174
$ws.Range("D4").Value = 174
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = $false
$ws.Range("H4").Value = 118
$ws.Range("J4").Value = "happy"
$ws.Range("K4").Value = 0.9615384615384616
$ws.Range("L4").Value = 25
$ws.Range("M4").Value = 25
$ws.Range("N4").Value = 1
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = $false
$ws.Range("Q4").Value = 1

# Row 5
$ws.Range("A5").Value = "panic"
$ws.Range("B5").Value = 0.1841085271317829
$ws.Range("C5").Value = 95
$ws.Range("D5").Value = 95
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = $false
$ws.Range("H5").Value = 421
$ws.Range("J5").Value = "best"
$ws.Range("K5").Value = 0.9322033898305084
$ws.Range("L5").Value = 55
$ws.Range("M5").Value = 55
$ws.Range("N5").Value = 1
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = $false
$ws.Range("Q5").Value = 4

# Row 6
$ws.Range("A6").Value = "sc"
$ws.Range("B6").Value = 0.1693121693121693
$ws.Range("C6").Value = 32
$ws.Range("D6").Value = 32
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = $false
$ws.Range("H6").Value = 157
$ws.Range("J6").Value = "love"
$ws.Range("K6").Value = 0.9130434782608695
$ws.Range("L6").Value = 42
$ws.Range("M6").Value = 42
$ws.Range("N6").Value = 1
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = $false
$ws.Range("Q6").Value = 4

# Row 7
$ws.Range("J7").Value = "great"
$ws.Range("K7").Value = 0.875
$ws.Range("L7").Value = 98
$ws.Range("M7").Value = 98
$ws.Range("N7").Value = 1
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = $false
$ws.Range("Q7").Value = 14

# Row 8
$ws.Range("J8").Value = "thanks"
$ws.Range("K8").Value = 0.8292682926829268
$ws.Range("L8").Value = 68
$ws.Range("M8").Value = 68
$ws.Range("N8").Value = 1
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = $false
$ws.Range("Q8").Value = 14

# Row 9
$ws.Range("J9").Value = "special"
$ws.Range("K9").Value = 0.8055555555555556
$ws.Range("L9").Value = 29
$ws.Range("M9").Value = 29
$ws.Range("N9").Value = 1
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = $false
$ws.Range("Q9").Value = 7

# Row 10
$ws.Range("J10").Value = "positive"
$ws.Range("K10").Value = 0.7931034482758621
$ws.Range("L10").Value = 46
$ws.Range("M10").Value = 46
$ws.Range("N10").Value = 1
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = $false
$ws.Range("Q10").Value = 12

# Row 11
$ws.Range("J11").Value = "thank"
$ws.Range("K11").Value = 0.7890625
$ws.Range("L11").Value = 101
$ws.Range("M11").Value = 101
$ws.Range("N11").Value = 1
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = $false
$ws.Range("Q11").Value = 27

# Row 12
$ws.Range("J12").Value = "free"
$ws.Range("K12").Value = 0.775
$ws.Range("L12").Value = 93
$ws.Range("M12").Value = 93
$ws.Range("N12").Value = 1
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = $false
$ws.Range("Q12").Value = 27

# Row 13
$ws.Range("J13").Value = "support"
$ws.Range("K13").Value = 0.7547169811320755
$ws.Range("L13").Value = 80
$ws.Range("M13").Value = 80
$ws.Range("N13").Value = 1
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = $false
$ws.Range("Q13").Value = 26

# Row 14
$ws.Range("J14").Value = "safe"
$ws.Range("K14").Value = 0.7535211267605634
$ws.Range("L14").Value = 107
$ws.Range("M14").Value = 107
$ws.Range("N14").Value = 1
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = $false
$ws.Range("Q14").Value = 35

# Row 15
$ws.Range("J15").Value = "safety"
$ws.Range("K15").Value = 0.7254901960784313
$ws.Range("L15").Value = 37
$ws.Range("M15").Value = 37
$ws.Range("N15").Value = 1
$ws.Range("O15").Value = 0
$ws.Range("P15").Value = $false
$ws.Range("Q15").Value = 14

# Row 16
$ws.Range("J16").Value = "confidence"
$ws.Range("K16").Value = 0.6944444444444444
$ws.Range("L16").Value = 25
$ws.Range("M16").Value = 25
$ws.Range("N16").Value = 1
$ws.Range("O16").Value = 0
$ws.Range("P16").Value = $false
$ws.Range("Q16").Value = 11

# Row 17
$ws.Range("J17").Value = "good"
$ws.Range("K17").Value = 0.68125
$ws.Range("L17").Value = 109
$ws.Range("M17").Value = 109
$ws.Range("N17").Value = 1
$ws.Range("O17").Value = 0
$ws.Range("P17").Value = $false
$ws.Range("Q17").Value = 51

# Row 18
$ws.Range("J18").Value = "better"
$ws.Range("K18").Value = 0.6190476190476191
$ws.Range("L18").Value = 39
$ws.Range("M18").Value = 39
$ws.Range("N18").Value = 1
$ws.Range("O18").Value = 0
$ws.Range("P18").Value = $false
$ws.Range("Q18").Value = 24

# Row 19
$ws.Range("J19").Value = "well"
$ws.Range("K19").Value = 0.6063829787234043
$ws.Range("L19").Value = 57
$ws.Range("M19").Value = 57
$ws.Range("N19").Value = 1
$ws.Range("O19").Value = 0
$ws.Range("P19").Value = $false
$ws.Range("Q19").Value = 37

# Row 20
$ws.Range("J20").Value = "relief"
$ws.Range("K20").Value = 0.58
$ws.Range("L20").Value = 29
$ws.Range("M20").Value = 29
$ws.Range("N20").Value = 1
$ws.Range("O20").Value = 0
$ws.Range("P20").Value = $false
$ws.Range("Q20").Value = 21

# Row 21
$ws.Range("J21").Value = "heroes"
$ws.Range("K21").Value = 0.5531914893617021
$ws.Range("L21").Value = 26
$ws.Range("M21").Value = 26
$ws.Range("N21").Value = 1
$ws.Range("O21").Value = 0
$ws.Range("P21").Value = $false
$ws.Range("Q21").Value = 21

# Row 22
$ws.Range("J22").Value = "hand"
$ws.Range("K22").Value = 0.5430809399477807
$ws.Range("L22").Value = 208
$ws.Range("M22").Value = 208
$ws.Range("N22").Value = 1
$ws.Range("O22").Value = 0
$ws.Range("P22").Value = $false
$ws.Range("Q22").Value = 175

# Row 23
$ws.Range("J23").Value = "fresh"
$ws.Range("K23").Value = 0.5416666666666666
$ws.Range("L23").Value = 26
$ws.Range("M23").Value = 26
$ws.Range("N23").Value = 1
$ws.Range("O23").Value = 0
$ws.Range("P23").Value = $false
$ws.Range("Q23").Value = 22

# Row 24
$ws.Range("J24").Value = "like"
$ws.Range("K24").Value = 0.4529411764705882
$ws.Range("L24").Value = 154
$ws.Range("M24").Value = 154
$ws.Range("N24").Value = 1
$ws.Range("O24").Value = 0
$ws.Range("P24").Value = $false
$ws.Range("Q24").Value = 186

# Row 25
$ws.Range("J25").Value = "care"
$ws.Range("K25").Value = 0.4157303370786517
$ws.Range("L25").Value = 37
$ws.Range("M25").Value = 37
$ws.Range("N25").Value = 1
$ws.Range("O25").Value = 0
$ws.Range("P25").Value = $false
$ws.Range("Q25").Value = 52

# Row 26
$ws.Range("J26").Value = "help"
$ws.Range("K26").Value = 0.4
$ws.Range("L26").Value = 118
$ws.Range("M26").Value = 118
$ws.Range("N26").Value = 1
$ws.Range("O26").Value = 0
$ws.Range("P26").Value = $false
$ws.Range("Q26").Value = 177

# Row 27
$ws.Range("J27").Value = "increase"
$ws.Range("K27").Value = 0.3717948717948718
$ws.Range("L27").Value = 29
$ws.Range("M27").Value = 29
$ws.Range("N27").Value = 1
$ws.Range("O27").Value = 0
$ws.Range("P27").Value = $false
$ws.Range("Q27").Value = 49

# Row 28
$ws.Range("J28").Value = "protect"
$ws.Range("K28").Value = 0.3561643835616438
$ws.Range("L28").Value = 26
$ws.Range("M28").Value = 26
$ws.Range("N28").Value = 1
$ws.Range("O28").Value = 0
$ws.Range("P28").Value = $false
$ws.Range("Q28").Value = 47

# Row 29
$ws.Range("J29").Value = "please"
$ws.Range("K29").Value = 0.3430962343096234
$ws.Range("L29").Value = 82
$ws.Range("M29").Value = 82
$ws.Range("N29").Value = 1
$ws.Range("O29").Value = 0
$ws.Range("P29").Value = $false
$ws.Range("Q29").Value = 157

# New row 29 needs the bold/bordered header style copied onto column J (matching J2:J28)
$ws.Range("J28").Copy()
$ws.Range("J29").PasteSpecial(-4122)
$excel.CutCopyMode = $false

Write-Output "edit complete"